$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N3").Value = 2.38
$ws.Range("O3").Value = 1.57
$ws.Range("S3").Value = 1.73
$ws.Range("S4").Value = 1.73
$ws.Range("J5").Value = 1.04
$ws.Range("K5").Value = 13
$ws.Range("L5").Value = 1.25
$ws.Range("M5").Value = 4
$ws.Range("N5").Value = 1.8
$ws.Range("P5").Value = 1.33
$ws.Range("Q5").Value = 3.25
$ws.Range("S5").Value = 1.91
$ws.Range("T5").Value = 7.5
$ws.Range("Y5").Value = 23
$ws.Range("Z5").Value = 13
$ws.Range("AA5").Value = 8
$ws.Range("AD5").Value = 251
$ws.Range("AF5").Value = 26
$ws.Range("N8").Value = 1.62
$ws.Range("G9").Value = 1.8
$ws.Range("X9").Value = 13
$ws.Range("G10").Value = 1.67
$ws.Range("G11").Value = 1.8
$ws.Range("J12").Value = 1.07
$ws.Range("K12").Value = 7.3
$ws.Range("L12").Value = 1.33
$ws.Range("M12").Value = 3.15
$ws.Range("N12").Value = 1.98
$ws.Range("O12").Value = 1.88
$ws.Range("P12").Value = 1.38
$ws.Range("Q12").Value = 2.87
$ws.Range("R12").Value = 1.72
$ws.Range("S12").Value = 2
$ws.Range("T12").Value = 8.5
$ws.Range("Y12").Value = 29
$ws.Range("Z12").Value = 10
$ws.Range("AD12").Value = 201
$ws.Range("AE12").Value = 9.5
$ws.Range("G13").Value = 2.22
$ws.Range("H13").Value = 3.05
$ws.Range("I13").Value = 3.1
$ws.Range("K13").Value = 6.3
$ws.Range("R13").Value = 1.87
$ws.Range("S13").Value = 1.83
$ws.Range("T13").Value = 6.9
$ws.Range("U13").Value = 10.25
$ws.Range("V13").Value = 9
$ws.Range("W13").Value = 22
$ws.Range("X13").Value = 19.5
$ws.Range("Z13").Value = 6.3
$ws.Range("AA13").Value = 6
$ws.Range("AB13").Value = 15
$ws.Range("AD13").Value = 700
$ws.Range("AE13").Value = 8.25
$ws.Range("AF13").Value = 15.5
$ws.Range("AG13").Value = 11.25
$ws.Range("AH13").Value = 40
$ws.Range("AI13").Value = 30
$ws.Range("AJ13").Value = 40
$ws.Range("G15").Value = 1.75
$ws.Range("H15").Value = 3.45
$ws.Range("I15").Value = 3.95
$ws.Range("L15").Value = 1.25
$ws.Range("M15").Value = 3.6
$ws.Range("N15").Value = 1.78
$ws.Range("O15").Value = 1.83
$ws.Range("R15").Value = 1.73
$ws.Range("S15").Value = 1.99
$ws.Range("T15").Value = 6.2
$ws.Range("U15").Value = 7.2
$ws.Range("V15").Value = 7
$ws.Range("W15").Value = 11.75
$ws.Range("X15").Value = 11.5
$ws.Range("Y15").Value = 19.5
$ws.Range("Z15").Value = 10.75
$ws.Range("AA15").Value = 6
$ws.Range("AB15").Value = 11.75
$ws.Range("AC15").Value = 45
$ws.Range("AD15").Value = 300
$ws.Range("AE15").Value = 10.25
$ws.Range("AF15").Value = 18.5
$ws.Range("AG15").Value = 11
$ws.Range("AH15").Value = 45
$ws.Range("AI15").Value = 28
$ws.Range("AJ15").Value = 30
$ws.Range("G16").Value = 1.09
$ws.Range("H16").Value = 7.2
$ws.Range("I16").Value = 20
$ws.Range("N16").Value = 1.29
$ws.Range("O16").Value = 3.3
$ws.Range("R16").Value = 2.31
$ws.Range("S16").Value = 1.54
$ws.Range("T16").Value = 9.5
$ws.Range("U16").Value = 6.1
$ws.Range("V16").Value = 10.25
$ws.Range("W16").Value = 5.6
$ws.Range("X16").Value = 9.25
$ws.Range("Y16").Value = 27
$ws.Range("Z16").Value = 22
$ws.Range("AA16").Value = 16
$ws.Range("AB16").Value = 30
$ws.Range("AC16").Value = 110
$ws.Range("AD16").Value = 600
$ws.Range("AE16").Value = 65
$ws.Range("AF16").Value = 250
$ws.Range("AG16").Value = 60
$ws.Range("AH16").Value = 101
$ws.Range("AI16").Value = 300
$ws.Range("AJ16").Value = 150
$ws.Range("I17").Value = 3.2
$ws.Range("J17").Value = 1.07
$ws.Range("K17").Value = 9
$ws.Range("L17").Value = 1.33
$ws.Range("M17").Value = 3.25
$ws.Range("N17").Value = 2.1
$ws.Range("O17").Value = 1.7
$ws.Range("P17").Value = 1.44
$ws.Range("Q17").Value = 2.63
$ws.Range("R17").Value = 1.8
$ws.Range("S17").Value = 1.95
$ws.Range("T17").Value = 7.5
$ws.Range("AB17").Value = 15
$ws.Range("AD17").Value = 251
$ws.Range("AE17").Value = 9.5
$ws.Range("G18").Value = 1.29
$ws.Range("R18").Value = 1.83
$ws.Range("S18").Value = 1.83
$ws.Range("G19").Value = 1.8
$ws.Range("L19").Value = 1.23
$ws.Range("R20").Value = 1.73
$ws.Range("J23").Value = 1.07
$ws.Range("K23").Value = 9
$ws.Range("L23").Value = 1.36
$ws.Range("M23").Value = 3
